# Update "template-soal" worksheet with a new question row (multi_choice)
# and move the active selection to A5, matching the authored workbook edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = "multi_choice"
$ws.Range("B4").Value = "Lorem ipsum?"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "hsdbfu"
$ws.Range("E4").Value = "yuegwruy"
$ws.Range("F4").Value = "uyweg"
$ws.Range("G4").Value = "weyrgwey"
$ws.Range("H4").Value = "a,b,d"

$ws.Range("A5").Select()
